$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 22.63
$ws.Range("P2").Value = 44.56
$ws.Range("Q2").Value = 67.19
$ws.Range("O3").Value = 23.26
$ws.Range("P3").Value = 44.36
$ws.Range("Q3").Value = 67.62
$ws.Range("O6").Value = 22.39
$ws.Range("P6").Value = 42.58
$ws.Range("Q6").Value = 64.97
$ws.Range("O7").Value = 22.3
$ws.Range("P7").Value = 43.52
$ws.Range("Q7").Value = 65.82
$ws.Range("O10").Value = 23.72
$ws.Range("P10").Value = 42.55
$ws.Range("Q10").Value = 66.27
$ws.Range("O11").Value = 24.16
$ws.Range("P11").Value = 43.34
$ws.Range("Q11").Value = 67.5
$ws.Range("O14").Value = 23.21
$ws.Range("P14").Value = 40.4
$ws.Range("Q14").Value = 63.61
$ws.Range("O15").Value = 23.58
$ws.Range("P15").Value = 40.65
$ws.Range("Q15").Value = 64.23
$ws.Range("O18").Value = 24.23
$ws.Range("P18").Value = 41.41
$ws.Range("Q18").Value = 65.64
$ws.Range("O19").Value = 24.33
$ws.Range("P19").Value = 42.52
$ws.Range("Q19").Value = 66.85
$ws.Range("O22").Value = 20.82
$ws.Range("P22").Value = 45.41
$ws.Range("Q22").Value = 66.23
$ws.Range("O23").Value = 20.69
$ws.Range("P23").Value = 45.67
$ws.Range("Q23").Value = 66.36
$ws.Range("O26").Value = 21.57
$ws.Range("P26").Value = 43.96
$ws.Range("Q26").Value = 65.53
$ws.Range("O27").Value = 21.75
$ws.Range("P27").Value = 44.23
$ws.Range("Q27").Value = 65.98
$ws.Range("O30").Value = 22.6
$ws.Range("P30").Value = 43.84
$ws.Range("Q30").Value = 66.44
$ws.Range("O31").Value = 22.61
$ws.Range("P31").Value = 42.48
$ws.Range("Q31").Value = 65.09
$ws.Range("O34").Value = 22.8
$ws.Range("P34").Value = 43.86
$ws.Range("Q34").Value = 66.66
$ws.Range("O35").Value = 23.24
$ws.Range("P35").Value = 41.41
$ws.Range("Q35").Value = 64.65
$ws.Range("O38").Value = 23.18
$ws.Range("P38").Value = 42.32
$ws.Range("Q38").Value = 65.5
$ws.Range("O39").Value = 22.89
$ws.Range("P39").Value = 43.77
$ws.Range("Q39").Value = 66.66
$ws.Range("O42").Value = 22.62
$ws.Range("P42").Value = 43.32
$ws.Range("Q42").Value = 65.94
$ws.Range("O43").Value = 22.58
$ws.Range("P43").Value = 44.7
$ws.Range("Q43").Value = 67.28
$ws.Range("O46").Value = 21.36
$ws.Range("P46").Value = 43.13
$ws.Range("Q46").Value = 64.49
$ws.Range("O47").Value = 22.57
$ws.Range("P47").Value = 42.53
$ws.Range("Q47").Value = 65.1
$ws.Range("O50").Value = 23.6
$ws.Range("P50").Value = 39.79
$ws.Range("Q50").Value = 63.39
$ws.Range("O51").Value = 23.13
$ws.Range("P51").Value = 41.62
$ws.Range("Q51").Value = 64.75
$ws.Range("O54").Value = 24.39
$ws.Range("P54").Value = 40.84
$ws.Range("Q54").Value = 65.23
$ws.Range("O55").Value = 24.51
$ws.Range("P55").Value = 41.96
$ws.Range("Q55").Value = 66.47
$ws.Range("O58").Value = 23.8
$ws.Range("P58").Value = 41.05
$ws.Range("Q58").Value = 64.85
$ws.Range("O59").Value = 24.04
$ws.Range("P59").Value = 40.24
$ws.Range("Q59").Value = 64.28
$ws.Range("O62").Value = 22.72
$ws.Range("P62").Value = 41.28
$ws.Range("Q62").Value = 64
$ws.Range("O63").Value = 22.66
$ws.Range("P63").Value = 41.24
$ws.Range("Q63").Value = 63.9
$ws.Range("O66").Value = 20.97
$ws.Range("P66").Value = 43.58
$ws.Range("Q66").Value = 64.55
$ws.Range("O67").Value = 21.19
$ws.Range("P67").Value = 43.18
$ws.Range("Q67").Value = 64.37
$ws.Range("O70").Value = 20.48
$ws.Range("P70").Value = 44.3
$ws.Range("Q70").Value = 64.78
$ws.Range("O71").Value = 21.27
$ws.Range("P71").Value = 43.72
$ws.Range("Q71").Value = 64.99
$ws.Range("O74").Value = 20.93
$ws.Range("P74").Value = 42.58
$ws.Range("Q74").Value = 63.51
$ws.Range("O75").Value = 21.17
$ws.Range("P75").Value = 44.18
$ws.Range("Q75").Value = 65.35
$ws.Range("O78").Value = 22.19
$ws.Range("P78").Value = 43.3
$ws.Range("Q78").Value = 65.49
$ws.Range("O79").Value = 21.66
$ws.Range("P79").Value = 42.44
$ws.Range("Q79").Value = 64.1
$ws.Range("O82").Value = 21.02
$ws.Range("P82").Value = 43.85
$ws.Range("Q82").Value = 64.87
$ws.Range("O83").Value = 21.54
$ws.Range("P83").Value = 41.64
$ws.Range("Q83").Value = 63.18
$ws.Range("O86").Value = 22.58
$ws.Range("P86").Value = 43.34
$ws.Range("Q86").Value = 65.92
$ws.Range("O87").Value = 22.86
$ws.Range("P87").Value = 43.92
$ws.Range("Q87").Value = 66.78
$ws.Range("O90").Value = 21.43
$ws.Range("P90").Value = 43.6
$ws.Range("Q90").Value = 65.03
$ws.Range("O91").Value = 22.26
$ws.Range("P91").Value = 42.19
$ws.Range("Q91").Value = 64.45
$ws.Range("O94").Value = 21.44
$ws.Range("P94").Value = 42.11
$ws.Range("Q94").Value = 63.55
$ws.Range("O95").Value = 22.74
$ws.Range("P95").Value = 40.92
$ws.Range("Q95").Value = 63.66
$ws.Range("O98").Value = 21.27
$ws.Range("P98").Value = 44.16
$ws.Range("Q98").Value = 65.43
$ws.Range("O99").Value = 20.63
$ws.Range("P99").Value = 45.05
$ws.Range("Q99").Value = 65.68
$ws.Range("O102").Value = 21.11
$ws.Range("P102").Value = 44.72
$ws.Range("Q102").Value = 65.83
$ws.Range("O103").Value = 22.28
$ws.Range("P103").Value = 43.44
$ws.Range("Q103").Value = 65.72
$ws.Range("O106").Value = 22.17
$ws.Range("P106").Value = 43.24
$ws.Range("Q106").Value = 65.41
$ws.Range("O107").Value = 22.46
$ws.Range("P107").Value = 43.54
$ws.Range("Q107").Value = 66
$ws.Range("O110").Value = 23.75
$ws.Range("P110").Value = 40.13
$ws.Range("Q110").Value = 63.88
$ws.Range("O111").Value = 23.02
$ws.Range("P111").Value = 40.13
$ws.Range("Q111").Value = 63.15
$ws.Range("O114").Value = 21.85
$ws.Range("P114").Value = 43.53
$ws.Range("Q114").Value = 65.38
$ws.Range("O115").Value = 22.44
$ws.Range("P115").Value = 42.32
$ws.Range("Q115").Value = 64.76
$ws.Range("O118").Value = 23.12
$ws.Range("P118").Value = 43.44
$ws.Range("Q118").Value = 66.56
$ws.Range("O119").Value = 22.75
$ws.Range("P119").Value = 43.69
$ws.Range("Q119").Value = 66.44
$ws.Range("O122").Value = 22.52
$ws.Range("P122").Value = 40.26
$ws.Range("Q122").Value = 62.78
$ws.Range("O123").Value = 22.86
$ws.Range("P123").Value = 40.7
$ws.Range("Q123").Value = 63.56
$ws.Range("O126").Value = 21.98
$ws.Range("P126").Value = 44.51
$ws.Range("Q126").Value = 66.49
$ws.Range("O127").Value = 22.4
$ws.Range("P127").Value = 42.3
$ws.Range("Q127").Value = 64.7
$ws.Range("O130").Value = 22.86
$ws.Range("P130").Value = 43.41
$ws.Range("Q130").Value = 66.27
$ws.Range("O131").Value = 22.77
$ws.Range("P131").Value = 42.54
$ws.Range("Q131").Value = 65.31
$ws.Range("O134").Value = 24.81
$ws.Range("P134").Value = 41.57
$ws.Range("Q134").Value = 66.38
$ws.Range("O135").Value = 24.45
$ws.Range("P135").Value = 42.11
$ws.Range("Q135").Value = 66.56
$ws.Range("O138").Value = 21.99
$ws.Range("P138").Value = 40.64
$ws.Range("Q138").Value = 62.63
$ws.Range("O139").Value = 22.47
$ws.Range("P139").Value = 39.71
$ws.Range("Q139").Value = 62.18
$ws.Range("O142").Value = 22.22
$ws.Range("P142").Value = 44.18
$ws.Range("Q142").Value = 66.4
$ws.Range("O143").Value = 23.55
$ws.Range("P143").Value = 43.12
$ws.Range("Q143").Value = 66.67
$ws.Range("O146").Value = 22.33
$ws.Range("P146").Value = 39.86
$ws.Range("Q146").Value = 62.19
$ws.Range("O147").Value = 23.07
$ws.Range("P147").Value = 39.82
$ws.Range("Q147").Value = 62.89
$ws.Range("O150").Value = 23.24
$ws.Range("P150").Value = 39.56
$ws.Range("Q150").Value = 62.8
$ws.Range("O151").Value = 24.72
$ws.Range("P151").Value = 38.61
$ws.Range("Q151").Value = 63.33
$ws.Range("O154").Value = 21.9
$ws.Range("P154").Value = 41.43
$ws.Range("Q154").Value = 63.33
$ws.Range("O155").Value = 21.77
$ws.Range("P155").Value = 42.16
$ws.Range("Q155").Value = 63.93
$ws.Range("O158").Value = 23.02
$ws.Range("P158").Value = 41.6
$ws.Range("Q158").Value = 64.62
$ws.Range("O159").Value = 22.96
$ws.Range("P159").Value = 41.84
$ws.Range("Q159").Value = 64.8
